# Fruta / hortaliza, semanal
# A new weekly price record is inserted as row 34; all subsequent rows
# (previously 34-73) shift down one position to become rows 35-74.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push rows 34..73 down to 35..74 and open up a blank row 34.
$ws.Rows.Item(34).Insert()

# Populate the newly inserted row 34 with the new weekly record.
$ws.Range("A34").Value = 3
$ws.Range("B34").Value = "Femacal de La Calera"
$ws.Range("C34").Value = "Coquimbo"
$ws.Range("D34").Value = 44895
$ws.Range("E34").Value = 5
$ws.Range("F34").Value = 100112022
$ws.Range("G34").Value = "Arveja Verde"
$ws.Range("H34").Value = "Perfection"
$ws.Range("I34").Value = "Primera"
$ws.Range("J34").Value = 73
$ws.Range("K34").Value = 22000
$ws.Range("L34").Value = 23000
$ws.Range("M34").Value = 22521
$ws.Range("N34").Value = "$/saco 25 kilos"
$ws.Range("O34").Value = "Región Metropolitana"
$ws.Range("P34").Value = 901
$ws.Range("Q34").Value = 25
$ws.Range("R34").Value = "Hortaliza"
